# Update the numeric dataframe values in columns G:R for rows 2-10
# (Cx_in, Cy_in, Cz_in, xA_in, yA_in, zA_in, xB_in, yB_in, zB_in, xC_in, yC_in, zC_in)
# to the recomputed values from the updated radical-dataframe generation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.002679, 0.116782, 0.419822, 1.4807, -0.700134, -0.125914, -1.425507, -0.784363, -0.123064, -0.057872, 1.367716, -0.170843),
    @(0.011203, 0.224431, 0.346664, -1.421002, -0.749111, -0.084967, 1.473656, -0.6730699999999999, -0.091923, -0.063856, 1.19775, -0.169773),
    @(-0.129309, 0.00706, 0.370343, 1.590831, -0.043715, -0.095932, -0.749094, -1.115596, -0.136415, -0.712429, 1.152252, -0.137996),
    @(-0.14432, 0.09843, 0.341551, 1.549684, -0.254352, -0.080599, -0.968772, -0.901166, -0.113244, -0.436591, 1.057088, -0.147709),
    @(-0.209128, 0.019171, 0.327533, 1.523963, -0.192741, -0.06741800000000001, -0.537409, 0.978013, -0.128213, -0.777426, -0.804442, -0.131902),
    @(0.00177, 0.009312000000000001, 0.353032, -0.662648, -1.117848, -0.115516, -0.637657, 1.127546, -0.119433, 1.298535, -0.019009, -0.118083),
    @(-0.004847, 0.08786099999999999, 0.333904, 1.16119, -0.533231, -0.099868, -1.084401, -0.647656, -0.099506, -0.071941, 1.093026, -0.13453),
    @(-0.09139600000000001, 0.004762, 0.312238, 1.227146, -0.047013, -0.0829, -0.5911920000000001, -0.884832, -0.114385, -0.544558, 0.927084, -0.114953),
    @(0.010979, -0.007273, 0.277261, 1.026126, -0.237963, -0.095317, -0.734105, -0.751586, -0.091922, -0.303, 0.996822, -0.090022)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $colNum = 7 + $c   # column G = 7
        $ws.Cells.Item($rowNum, $colNum).Value = $rowValues[$c]
    }
}
